$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.897.23"
$ws.Range("E2").Value = "'  -0.54%  "

$ws.Range("D3").Value = "'1.906.16"
$ws.Range("E3").Value = "'  -0.29%  "

$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = "'  -0.60%  "

$ws.Range("D5").Value = "'313.38"
$ws.Range("E5").Value = "'  -0.73%  "

$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = "'  -0.53%  "

$ws.Range("D7").Value = "'0.4985"
$ws.Range("E7").Value = "'  +3.54%  "

$ws.Range("E8").Value = "'  -0.14%  "

$ws.Range("D9").Value = "'0.07288"
$ws.Range("E9").Value = "'  -1.00%  "

$ws.Range("D10").Value = "'0.9121"
$ws.Range("E10").Value = "'  -2.37%  "

$ws.Range("D11").Value = "'21.10"
$ws.Range("E11").Value = "'  +1.16%  "

$ws.Range("B12").Value = "'TRON"
$ws.Range("C12").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07686"
$ws.Range("E12").Value = "'  -1.25%  "

$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.922.09"
$ws.Range("E13").Value = "'  +0.45%  "

$ws.Range("D14").Value = "'5.499"
$ws.Range("E14").Value = "'  -0.01%  "

$ws.Range("D15").Value = "'92.22"
$ws.Range("E15").Value = "'  +0.06%  "

$ws.Range("D16").Value = "'0.9992"
$ws.Range("E16").Value = "'  -0.70%  "

$ws.Range("D17").Value = "'0.000008730"
$ws.Range("E17").Value = "'  -1.56%  "

$ws.Range("D18").Value = "'0.9990"
$ws.Range("E18").Value = "'  -0.51%  "

$ws.Range("D19").Value = "'27.929.36"
$ws.Range("E19").Value = "'  -0.57%  "

$ws.Range("D20").Value = "'14.64"
$ws.Range("E20").Value = "'  -0.97%  "

$ws.Range("D21").Value = "'5.179"
$ws.Range("E21").Value = "'  +0.30%  "

$ws.Range("E22").Value = "'  -0.67%  "

$ws.Range("D23").Value = "'6.571"
$ws.Range("E23").Value = "'  -1.12%  "

$ws.Range("D24").Value = "'153.00"
$ws.Range("E24").Value = "'  -1.86%  "

$ws.Range("D25").Value = "'1.863"
$ws.Range("E25").Value = "'  -2.92%  "

$ws.Range("E26").Value = "'  +4.16%  "

$ws.Range("E27").Value = "'  -0.47%  "

$ws.Range("D28").Value = "'115.35"
$ws.Range("E28").Value = "'  -1.32%  "

$ws.Range("D29").Value = "'4.900"
$ws.Range("E29").Value = "'  -1.39%  "

$ws.Range("D30").Value = "'0.09010"
$ws.Range("E30").Value = "'  +0.57%  "

$ws.Range("D31").Value = "'3.201"

$ws.Range("D32").Value = "'4.860"
$ws.Range("E32").Value = "'  +3.87%  "

$ws.Range("D33").Value = "'1.233"
$ws.Range("E33").Value = "'  -2.75%  "

$ws.Range("D34").Value = "'0.7740"
$ws.Range("E34").Value = "'  -0.87%  "

$ws.Range("D35").Value = "'0.02089"
$ws.Range("E35").Value = "'  +1.57%  "

$ws.Range("B36").Value = "'RenderToken"
$ws.Range("C36").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'2.560"
$ws.Range("E36").Value = "'  -2.24%  "

$ws.Range("B37").Value = "'MXToken"
$ws.Range("C37").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'3.062"
$ws.Range("E37").Value = "'  +2.24%  "

$ws.Range("B38").Value = "'TrustWalletToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.094"
$ws.Range("E38").Value = "'  -1.73%  "

$ws.Range("B39").Value = "'TheSandbox"
$ws.Range("C39").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.5551"
$ws.Range("E39").Value = "'  +0.83%  "

$ws.Range("B40").Value = "'Hedera"
$ws.Range("C40").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "'0.05280"
$ws.Range("E40").Value = "'  -0.72%  "

$ws.Range("B41").Value = "'FraxShare"
$ws.Range("C41").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.889"
$ws.Range("E41").Value = "'  -2.05%  "

$ws.Range("B42").Value = "'Aptos"
$ws.Range("C42").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'8.507"
$ws.Range("E42").Value = "'  +0.16%  "

$ws.Range("B43").Value = "'Algorand"
$ws.Range("C43").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1521"
$ws.Range("E43").Value = "'  -0.53%  "

$ws.Range("B44").Value = "'Quant"
$ws.Range("C44").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'112.66"
$ws.Range("E44").Value = "'  +3.87%  "

$ws.Range("B45").Value = "'EnergySwap"
$ws.Range("C45").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'10.62"
$ws.Range("E45").Value = "'  -0.79%  "

$ws.Range("D46").Value = "'0.4840"
$ws.Range("E46").Value = "'  +0.20%  "

$ws.Range("B47").Value = "'PaxDollar"
$ws.Range("C47").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'0.9989"
$ws.Range("E47").Value = "'  -0.57%  "

$ws.Range("B48").Value = "'NEARProtocol"
$ws.Range("C48").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.636"
$ws.Range("E48").Value = "'  -0.89%  "

$ws.Range("B49").Value = "'Aave"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'67.50"
$ws.Range("E49").Value = "'  -0.83%  "

$ws.Range("B50").Value = "'Cronos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06052"
$ws.Range("E50").Value = "'  -0.49%  "

$ws.Range("B51").Value = "'EOS"
$ws.Range("C51").Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'0.9067"
$ws.Range("E51").Value = "'  +0.74%  "
